$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply each cell update. Columns D (Price) hold numeric-looking text that must
# be preserved exactly as text (e.g. "1.00", "11.10", "3.51"), and some D values
# use multi-dot "thousands" notation (e.g. "56.557.44") which already stays text.
# Column B/C are plain text, column E percentage strings already stay text because
# of the surrounding spaces and "%" suffix. To be safe and consistent we force all
# D-column price cells to Text format before assignment, then reset the style back
# to Normal so no stray cell formatting/style indices are left behind.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '56.557.44'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +10.50%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.257.42'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +6.30%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '399.57'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.12%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '111.15'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +9.41%  '
$ws.Range("E7").Value = '  +4.44%  '
$ws.Range("E8").Value = '  -0.05%  '
$ws.Range("E9").Value = '  +6.77%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '39.55'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +7.68%  '
$ws.Range("E11").Value = '  +11.75%  '
$ws.Range("E12").Value = '  +2.36%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.771.69'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +6.43%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '19.26'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +5.33%  '
$ws.Range("E15").Value = '  +5.80%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.258.41'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +6.44%  '
$ws.Range("E17").Value = '  +5.45%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '11.01'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +3.93%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '56.420.58'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +10.28%  '
$ws.Range("E20").Value = '  +4.89%  '
$ws.Range("E21").Value = '  +8.62%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '13.06'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +6.57%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '297.71'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +12.77%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '75.16'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +7.88%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.22'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.65%  '
$ws.Range("E26").Value = '  +3.09%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '28.26'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +5.58%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '4.36'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +4.66%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.35'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.63%  '
$ws.Range("E30").Value = '  +3.84%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.999'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.09%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.111'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +6.43%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '11.10'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +6.00%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '38.16'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +6.51%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0489'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.29%  '
$ws.Range("E36").Value = '  +3.56%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '51.79'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +3.76%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.13'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +26.11%  '
$ws.Range("B39").Value = 'FirstDigitalUSD'
$ws.Range("C39").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.00'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.03%  '
$ws.Range("B40").Value = 'LidoDAOToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.51'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +4.03%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '17.62'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +6.48%  '
$ws.Range("E42").Value = '  +6.35%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '133.31'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.86%  '
$ws.Range("E44").Value = '  +5.30%  '
$ws.Range("E45").Value = '  +4.66%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.285'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.40%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '22.22'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.37%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.13'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +50.06%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.149.26'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +4.07%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.09'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.54%  '
$ws.Range("E51").Value = '  -2.67%  '
